# Weekly Fruit/Vegetable price update: a new weekly price record is
# inserted into the "Arveja Verde" (Vega Central Mapocho de Santiago) data
# set. The new observation belongs right after the existing row for date
# serial 44483 (row 90) and before the former row 90 (date 44489), so a
# whole new row is inserted at position 90 and every subsequent record
# shifts down by one. The sheet's used range therefore grows from
# A1:R107 to A1:R108.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 90, pushing the existing rows 90-107
# down to 91-108.
$ws.Rows.Item(90).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A90").Value = 9
$ws.Range("B90").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C90").Value = "Metropolitana"
$ws.Range("D90").Value = 44617
$ws.Range("E90").Value = 13
$ws.Range("F90").Value = 100112022
$ws.Range("G90").Value = "Arveja Verde"
$ws.Range("H90").Value = "Sin especificar"
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 28
$ws.Range("K90").Value = 25000
$ws.Range("L90").Value = 26000
$ws.Range("M90").Value = 25500
$ws.Range("N90").Value = "$/saco 25 kilos"
$ws.Range("O90").Value = "Carahue"
$ws.Range("P90").Value = 1020
$ws.Range("Q90").Value = 25
$ws.Range("R90").Value = "Hortaliza"
